# Apply the "missing_data" re-randomization edit described by the diff.
#
# Summary of the change:
#  1. Two data rows are removed entirely from the sheet: "RM 232" (row 26)
#     and "SC 92" (row 28). All rows below them shift up by the
#     corresponding amount, and the sheet dimension shrinks from
#     A1:F35 to A1:F33.
#  2. A handful of individual cells in columns D/E (and a couple in B)
#     switch between "missing" (blank) and a concrete numeric value,
#     reflecting a different random "missingness" pattern for the
#     remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the two obsolete rows -------------------------------
# Delete the higher-numbered row first so the lower row index ("RM 232",
# row 26) stays valid while we work.
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# --- Step 2: update individual cells to match the new missing pattern ---
# Helper references: columns B=2, D=4, E=5

# Row 2
$ws.Cells.Item(2, 4).Value = -13.5

# Row 3
$ws.Cells.Item(3, 4).Value = ""

# Row 4
$ws.Cells.Item(4, 4).Value = ""

# Row 5
$ws.Cells.Item(5, 5).Value = ""

# Row 8
$ws.Cells.Item(8, 5).Value = -6.6

# Row 10
$ws.Cells.Item(10, 5).Value = -6.1

# Row 11
$ws.Cells.Item(11, 4).Value = -15.5

# Row 12
$ws.Cells.Item(12, 5).Value = ""

# Row 13
$ws.Cells.Item(13, 4).Value = ""

# Row 15
$ws.Cells.Item(15, 5).Value = -8.4

# Row 18
$ws.Cells.Item(18, 5).Value = ""

# Row 19
$ws.Cells.Item(19, 5).Value = ""

# Row 21
$ws.Cells.Item(21, 4).Value = -14.3

# Row 25
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = -7.1

# Row 27 ("SC 101" after the row deletions)
$ws.Cells.Item(27, 5).Value = -10

# Row 29 ("SC 119" after the row deletions)
$ws.Cells.Item(29, 2).Value = ""
$ws.Cells.Item(29, 5).Value = ""

# Row 33 ("SC 232" after the row deletions)
$ws.Cells.Item(33, 2).Value = -19.5
$ws.Cells.Item(33, 4).Value = -14.1
$ws.Cells.Item(33, 5).Value = ""

Write-Output "Edit applied."
